$wb = $excel.ActiveWorkbook

# 1. Update the status text "Ready for handoff" -> "In Translation" on every sheet.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $used.Replace("Ready for handoff", "In Translation")
}

# 2. The status columns are narrower now that the longest status string is
#    shorter ("In Translation" vs. "Ready for handoff"), so re-fit their width.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E:F").ColumnWidth = 12.5

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C:C").ColumnWidth = 12.5

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C:C").ColumnWidth = 12.5
